# Location geändert, zwei Beispiel-TaxiStands eingefügt
$wb = $excel.ActiveWorkbook

# --- Sheet "TaxiStands": update the two existing taxi stand coordinates ---
$ws2 = $wb.Worksheets.Item("TaxiStands")
$ws2.Cells.Item(2, 1).Value = 52.381989108834198
$ws2.Cells.Item(2, 2).Value = 9.7381182223916092
$ws2.Cells.Item(3, 1).Value = 52.375468581529297
$ws2.Cells.Item(3, 2).Value = 9.7509848120552896

# --- Sheet "Taxis": update the taxi's current location ---
$ws3 = $wb.Worksheets.Item("Taxis")
$ws3.Cells.Item(2, 1).Value = 52.381989108834198
$ws3.Cells.Item(2, 2).Value = 9.7401182223916098

# --- Sheet "Orders": move all orders to the new location ---
$ws1 = $wb.Worksheets.Item("Orders")
for ($r = 2; $r -le 9; $r++) {
    $ws1.Cells.Item($r, 1).Value = 120
    $ws1.Cells.Item($r, 2).Value = 52.375394479042797
    $ws1.Cells.Item($r, 3).Value = 9.7315180260351593
    $ws1.Cells.Item($r, 4).Value = 52.382591097574597
    $ws1.Cells.Item($r, 5).Value = 9.7309718027690906
}

# --- Restore the view / selection state on each sheet ---
$ws2.Activate() | Out-Null
$ws2.Range("C6").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("B2").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("A9").Select() | Out-Null
